$d = $word.ActiveDocument

# The requirement sentence that gains the word "ver" and is wrapped in the
# "_Hlk131584703" bookmark (the auto-generated bookmark Word drops in when
# text that was copied/pasted gets referenced later).
$sentence = "El sistema debe permitir al dueño las facturas y seleccionar alguna y generarla nuevamente, así mismo enviarla"

$sentenceRange = $d.Content
$found = $sentenceRange.Find.Execute($sentence, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, "", 0)

if ($found) {
    $startPos = $sentenceRange.Start
    $endPos = $sentenceRange.End

    # Wrap the whole sentence in the bookmark first; Word keeps bookmark
    # ranges in sync with edits made inside them, so it will grow to cover
    # the new word we are about to insert.
    $bookmarkRange = $d.Range($startPos, $endPos)
    $d.Bookmarks.Add("_Hlk131584703", $bookmarkRange) | Out-Null

    # Within that same sentence, drop the cursor right after "dueño" (before
    # its trailing space) and type " ver", turning "...al dueño las
    # facturas..." into "...al dueño ver las facturas...".
    $wordRange = $d.Range($startPos, $endPos)
    $wordRange.Find.Execute("dueño", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0) | Out-Null
    $insertionPoint = $d.Range($wordRange.End, $wordRange.End)
    $insertionPoint.InsertAfter(" ver")
}
